{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find and delete the empty paragraph (used to generate an empty line\n// for an empty AQL expression result).\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  if (paragraphs.items[i].text === \"\") {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the empty paragraph (was used to render an empty line for an\n# empty AQL expression result).\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq [char]13) {\n        $p.Range.Delete()\n    }\n}\n"}
